$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new record at row 177 (shifts old rows 177..289 down to 178..290) ---
$ws.Rows.Item(177).Insert()

$ws.Range("A177").Value = 10
$ws.Range("B177").Value = 'Vega Modelo de Temuco'
$ws.Range("C177").Value = 'La Araucanía'
$ws.Range("D177").Value = 44567
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = 100112040
$ws.Range("G177").Value = 'Cilantro'
$ws.Range("H177").Value = 'Sin especificar'
$ws.Range("I177").Value = 'Primera'
$ws.Range("J177").Value = 70
$ws.Range("K177").Value = 8000
$ws.Range("L177").Value = 8000
$ws.Range("M177").Value = 8000
$ws.Range("N177").Value = '$/docena de atados (2 kilos)'
$ws.Range("O177").Value = 'Provincia de Cautín'
$ws.Range("P177").Value = 4000
$ws.Range("Q177").Value = 2
$ws.Range("R177").Value = 'Hortaliza'

# --- Insert second new record at row 281 (shifts rows currently 281..290 down to 282..291) ---
$ws.Rows.Item(281).Insert()

$ws.Range("A281").Value = 10
$ws.Range("B281").Value = 'Vega Modelo de Temuco'
$ws.Range("C281").Value = 'La Araucanía'
$ws.Range("D281").Value = 44568
$ws.Range("E281").Value = 9
$ws.Range("F281").Value = 100112040
$ws.Range("G281").Value = 'Cilantro'
$ws.Range("H281").Value = 'Sin especificar'
$ws.Range("I281").Value = 'Primera'
$ws.Range("J281").Value = 30
$ws.Range("K281").Value = 8000
$ws.Range("L281").Value = 8000
$ws.Range("M281").Value = 8000
$ws.Range("N281").Value = '$/docena de atados (2 kilos)'
$ws.Range("O281").Value = 'Provincia de Cautín'
$ws.Range("P281").Value = 4000
$ws.Range("Q281").Value = 2
$ws.Range("R281").Value = 'Hortaliza'

# Ensure the date-formatted column D keeps consistent number formatting with the rest of column D
$ws.Range("D177").NumberFormat = $ws.Range("D178").NumberFormat
$ws.Range("D281").NumberFormat = $ws.Range("D280").NumberFormat
